$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '41.511.28'
$ws.Range('E2').Value = '  -2.10%  '
$ws.Range('D3').Value = '2.483.49'
$ws.Range('E3').Value = '  -0.99%  '
$ws.Range('E4').Value = '  +0.18%  '
$ws.Range('D5').Value = '314.03'
$ws.Range('E5').Value = '  +1.24%  '
$ws.Range('D6').Value = '94.29'
$ws.Range('E6').Value = '  -3.91%  '
$ws.Range('E7').Value = '  -1.86%  '
$ws.Range('E8').Value = '  +0.22%  '
$ws.Range('D9').Value = '0.499'
$ws.Range('E9').Value = '  -2.75%  '
$ws.Range('D10').Value = '33.50'
$ws.Range('E10').Value = '  -4.16%  '
$ws.Range('D11').Value = '0.0780'
$ws.Range('E11').Value = '  -1.87%  '
$ws.Range('E12').Value = '  +0.28%  '
$ws.Range('B13').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C13').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D13').Value = '2.866.49'
$ws.Range('E13').Value = '  -0.79%  '
$ws.Range('B14').Value = 'Polkadot'
$ws.Range('C14').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D14').Value = '6.91'
$ws.Range('E14').Value = '  -3.15%  '
$ws.Range('D15').Value = '15.43'
$ws.Range('E15').Value = '  +0.41%  '
$ws.Range('D16').Value = '2.417.76'
$ws.Range('E16').Value = '  -3.62%  '
$ws.Range('D17').Value = '0.792'
$ws.Range('E17').Value = '  -1.37%  '
$ws.Range('D18').Value = '41.455.21'
$ws.Range('E18').Value = '  -2.21%  '
$ws.Range('D19').Value = '6.33'
$ws.Range('E20').Value = '  -1.14%  '
$ws.Range('D21').Value = '11.24'
$ws.Range('E21').Value = '  -5.10%  '
$ws.Range('D22').Value = '68.91'
$ws.Range('E22').Value = '  +0.65%  '
$ws.Range('D23').Value = '236.82'
$ws.Range('E23').Value = '  -1.39%  '
$ws.Range('D24').Value = '2.76'
$ws.Range('E24').Value = '  -2.29%  '
$ws.Range('B25').Value = 'Dai'
$ws.Range('C25').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D25').Value = '1.00'
$ws.Range('E25').Value = '  +0.06%  '
$ws.Range('B26').Value = 'ImmutableX'
$ws.Range('C26').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D26').Value = '1.90'
$ws.Range('E26').Value = '  -3.20%  '
$ws.Range('D27').Value = '24.11'
$ws.Range('E27').Value = '  -3.78%  '
$ws.Range('D28').Value = '2.25'
$ws.Range('E28').Value = '  -0.18%  '
$ws.Range('D29').Value = '9.78'
$ws.Range('E29').Value = '  -1.72%  '
$ws.Range('D30').Value = '36.63'
$ws.Range('E30').Value = '  -2.92%  '
$ws.Range('D31').Value = '152.65'
$ws.Range('E31').Value = '  -2.69%  '
$ws.Range('D32').Value = '5.48'
$ws.Range('E32').Value = '  -5.31%  '
$ws.Range('D33').Value = '2.57'
$ws.Range('E33').Value = '  -2.97%  '
$ws.Range('D34').Value = '18.22'
$ws.Range('E34').Value = '  +5.83%  '
$ws.Range('D35').Value = '0.0756'
$ws.Range('E35').Value = '  -2.82%  '
$ws.Range('D36').Value = '3.08'
$ws.Range('E36').Value = '  -0.88%  '
$ws.Range('D37').Value = '2.47'
$ws.Range('E37').Value = '  -12.58%  '
$ws.Range('D38').Value = '1.88'
$ws.Range('E38').Value = '  -2.91%  '
$ws.Range('E39').Value = '  -1.74%  '
$ws.Range('E40').Value = '  -4.51%  '
$ws.Range('D41').Value = '4.16'
$ws.Range('E41').Value = '  +0.66%  '
$ws.Range('D42').Value = '1.01'
$ws.Range('E42').Value = '  +0.27%  '
$ws.Range('D43').Value = '19.70'
$ws.Range('E43').Value = '  -6.95%  '
$ws.Range('D44').Value = '1.989.03'
$ws.Range('E44').Value = '  -0.44%  '
$ws.Range('D45').Value = '0.0286'
$ws.Range('E45').Value = '  -1.96%  '
$ws.Range('D46').Value = '3.02'
$ws.Range('E46').Value = '  -6.10%  '
$ws.Range('D47').Value = '8.89'
$ws.Range('E47').Value = '  -2.03%  '
$ws.Range('D48').Value = '2.730.53'
$ws.Range('E48').Value = '  -0.81%  '
$ws.Range('D49').Value = '69.56'
$ws.Range('E49').Value = '  -1.47%  '
$ws.Range('D50').Value = '97.21'
$ws.Range('E50').Value = '  -2.35%  '
$ws.Range('D51').Value = '0.178'
$ws.Range('E51').Value = '  -4.44%  '
